# Weekly update: insert a new week's record at row 76 for
# "Hortaliza, Terminal La Palmera de La Serena - Perejil".
# This pushes the existing rows 76-135 down to 77-136 (dates/prices
# keep their original relative order - each old row's data reappears
# one row lower), and the freed-up row 76 gets the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 76, shifting rows 76:135 down
# to 77:136 (matches native Excel "Insert Sheet Rows" behaviour).
$ws.Rows.Item(76).Insert()

# Populate the newly inserted row 76 with the new record.
$ws.Cells.Item(76, 1).Value = 8
$ws.Cells.Item(76, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(76, 3).Value = "Coquimbo"
$ws.Cells.Item(76, 4).Value = 44658
$ws.Cells.Item(76, 5).Value = 4
$ws.Cells.Item(76, 6).Value = 100112044
$ws.Cells.Item(76, 7).Value = "Perejil"
$ws.Cells.Item(76, 8).Value = "Sin especificar"
$ws.Cells.Item(76, 9).Value = "Primera"
$ws.Cells.Item(76, 10).Value = 3000
$ws.Cells.Item(76, 11).Value = 2000
$ws.Cells.Item(76, 12).Value = 2500
$ws.Cells.Item(76, 13).Value = 2250
$ws.Cells.Item(76, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(76, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(76, 16).Value = 1500
$ws.Cells.Item(76, 17).Value = 1.5
$ws.Cells.Item(76, 18).Value = "Hortaliza"

Write-Output "Inserted new row 76; sheet now spans $($ws.UsedRange.Rows.Count) rows"
